$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shorten the "Previous land use" wording (col F) and fix a typo in A3,
# per the manuscript revision for Grace.
$ws.Range("A3").Value = "Silvopasture (SP)"
$ws.Range("F2").Value = "SP (pasture, savanna or secondary forest), FG"
$ws.Range("F4").Value = "CF,  SP (pasture, savanna or secondary forest), FG"
$ws.Range("F5").Value = "CF, YA"
$ws.Range("F6").Value = "HG, YA"

# Leave the view scrolled back to the top with G6 as the active selection.
$ws.Activate()
$ws.Range("G6").Select()
